$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.746142864227295
$ws.Range("B1").Value = 2.843096971511841
$ws.Range("C1").Value = 2.490867853164673
$ws.Range("D1").Value = 1.664771676063538
$ws.Range("E1").Value = 0.793475866317749
